# "no matter handle visiable or not, slider will be supported"
#
# ScrollView (row 6), Scrollbar (row 7) and Slider (row 8) are marked as
# done ("V" in the 【完成状态】 column), and the Scrollbar's required-child
# prefixes drop the now-unneeded "f_" entry. A new reference row is added
# at the bottom of the sheet documenting the PSDName.xml / Globle naming
# convention, highlighted with Excel's built-in "Good" (green) cell style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark ScrollView / Scrollbar / Slider as complete.
$ws.Range("F6").Value = "V"
$ws.Range("F7").Value = "V"
$ws.Range("F8").Value = "V"

# Scrollbar's required-sub-object prefixes: b_,f_,h_ -> b_,h_
$ws.Range("D7").Value = "b_,h_"

# New documentation row explaining PSDName.xml / Globle-named xml behaviour.
$ws.Range("A25").Value = "PSDName.xml"
$ws.Range("B25").Value = "带Globle名的xml，会将图片导入到指定文件夹中"
$ws.Range("B25").Style = "Good"

# Move the active selection the way the author left it (F7 -> F8).
$ws.Range("F8").Select() | Out-Null
